# Register by email, small ui fixes

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Two more conversations collected ("register by email" flow).
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "conversation_11_07_2023__09_31_44"
$ws.Range("C10").Value = "ai"

$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "conversation_11_07_2023__09_38_05"
$ws.Range("C11").Value = "ai"

# Small UI fixes: widen the conversation-name column so names aren't
# truncated, and leave the selection on the column with the new data.
$ws.Columns.Item(2).ColumnWidth = 43.71
$ws.Range("B16").Select()

# Best-effort: match the maximized/full-screen app window from the source edit.
$win = $excel.ActiveWindow
$win.WindowState = -4143
$win.Left = -110
$win.Top = -110
$win.Width = 1420
$win.Height = 770
